# Auto-generated edit script applying numeric corrections to Mateus_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 51
  $ws.Range("H51").Value = 11735.777
  $ws.Range("J51").Value = 14145.75
  $ws.Range("L51").Value = 14145.75
  $ws.Range("N51").Value = -15113.75
  # Row 116
  $ws.Range("H116").Value = 4966.6665
  $ws.Range("J116").Value = 4966.6665
  $ws.Range("L116").Value = 4966.6665
  $ws.Range("N116").Value = -11850.6665
  # Row 132
  $ws.Range("H132").Value = 11981.85
  $ws.Range("I132").Value = 8329.134
  $ws.Range("J132").Value = 22940
  $ws.Range("K132").Value = 24987.402
  $ws.Range("L132").Value = 68820
  $ws.Range("M132").Value = -22457.402
  $ws.Range("N132").Value = -73880
  # Row 137
  $ws.Range("H137").Value = 2302.9119
  $ws.Range("J137").Value = 2493.25
  $ws.Range("L137").Value = 7479.75
  $ws.Range("N137").Value = -12579.75
  # Row 138
  $ws.Range("H138").Value = 3176.375
  $ws.Range("I138").Value = 4097
  $ws.Range("J138").Value = 3115
  $ws.Range("K138").Value = 12291
  $ws.Range("L138").Value = 9345
  $ws.Range("M138").Value = -7151
  $ws.Range("N138").Value = -19625
  # Row 141
  $ws.Range("H141").Value = 1615.6
  $ws.Range("I141").Value = 1640.8485
  $ws.Range("K141").Value = 4922.5455
  $ws.Range("M141").Value = 257.4544999999998

$ws = $wb.Worksheets.Item("ARM")
  # Row 32
  $ws.Range("H32").Value = 4048.2056
  $ws.Range("I32").Value = 3959.6057
  $ws.Range("K32").Value = 3959.6057
  $ws.Range("M32").Value = -3672.6057
  # Row 41
  $ws.Range("H41").Value = 27975
  $ws.Range("I41").Value = 2850
  $ws.Range("J41").Value = 33000
  $ws.Range("K41").Value = 2850
  $ws.Range("L41").Value = 33000
  $ws.Range("M41").Value = -2436
  $ws.Range("N41").Value = -33828
  # Row 45
  $ws.Range("H45").Value = 99403
  $ws.Range("I45").Value = 288782.72
  $ws.Range("J45").Value = 4713.143
  $ws.Range("K45").Value = 288782.72
  $ws.Range("L45").Value = 4713.143
  $ws.Range("M45").Value = -288405.72
  $ws.Range("N45").Value = -5467.143
  # Row 61
  $ws.Range("H61").Value = 9649.6
  $ws.Range("I61").Value = 5617.206
  $ws.Range("K61").Value = 5617.206
  $ws.Range("M61").Value = -5405.206
  # Row 122
  $ws.Range("H122").Value = 2132.8572
  $ws.Range("I122").Value = 1786.8
  $ws.Range("K122").Value = 5360.4
  $ws.Range("M122").Value = -2910.4
  # Row 132
  $ws.Range("H132").Value = 3958.48
  $ws.Range("J132").Value = 2740
  $ws.Range("L132").Value = 8220
  $ws.Range("N132").Value = -13280
  # Row 136
  $ws.Range("H136").Value = 9649.6
  $ws.Range("I136").Value = 5617.206
  $ws.Range("K136").Value = 16851.618
  $ws.Range("M136").Value = -14301.618

$ws = $wb.Worksheets.Item("BSM")
  # Row 2
  $ws.Range("H2").Value = 62806.332
  $ws.Range("J2").Value = 69209.5
  $ws.Range("L2").Value = 69209.5
  $ws.Range("N2").Value = -69435.5
  # Row 20
  $ws.Range("H20").Value = 4108.8823
  $ws.Range("I20").Value = 3334.7693
  $ws.Range("K20").Value = 3334.7693
  $ws.Range("M20").Value = -3087.7693
  # Row 92
  $ws.Range("H92").Value = 61851.668
  $ws.Range("J92").Value = 61851.668
  $ws.Range("L92").Value = 61851.668
  $ws.Range("N92").Value = -66843.66800000001
  # Row 110
  $ws.Range("H110").Value = 0
  $ws.Range("J110").Value = 0
  $ws.Range("L110").Value = 0
  $ws.Range("N110").Value = $null
  # Row 134
  $ws.Range("H134").Value = 3262.5715
  $ws.Range("I134").Value = 3322.4102
  $ws.Range("J134").Value = 2484.6667
  $ws.Range("K134").Value = 9967.230599999999
  $ws.Range("L134").Value = 7454.000100000001
  $ws.Range("M134").Value = -7432.230599999999
  $ws.Range("N134").Value = -12524.0001

$ws = $wb.Worksheets.Item("CRP")
  # Row 16
  $ws.Range("H16").Value = 4883.25
  $ws.Range("J16").Value = 5673.625
  $ws.Range("L16").Value = 5673.625
  $ws.Range("N16").Value = -6247.625
  # Row 31
  $ws.Range("H31").Value = 6376.8887
  $ws.Range("I31").Value = 5267.2856
  $ws.Range("K31").Value = 5267.2856
  $ws.Range("M31").Value = -4972.2856
  # Row 34
  $ws.Range("H34").Value = 6376.8887
  $ws.Range("I34").Value = 5267.2856
  $ws.Range("K34").Value = 5267.2856
  $ws.Range("M34").Value = -5065.2856
  # Row 58
  $ws.Range("H58").Value = 5273.7144
  $ws.Range("I58").Value = 4072.5
  $ws.Range("K58").Value = 4072.5
  $ws.Range("M58").Value = -3869.5
  # Row 74
  $ws.Range("H74").Value = 42524.5
  $ws.Range("J74").Value = 42524.5
  $ws.Range("L74").Value = 42524.5
  $ws.Range("N74").Value = -44272.5
  # Row 77
  $ws.Range("H77").Value = 42524.5
  $ws.Range("J77").Value = 42524.5
  $ws.Range("L77").Value = 127573.5
  $ws.Range("N77").Value = -136309.5
  # Row 105
  $ws.Range("H105").Value = 1972.6666
  $ws.Range("I105").Value = 1407.8572
  $ws.Range("K105").Value = 1407.8572
  $ws.Range("M105").Value = 339.1428000000001
  # Row 112
  $ws.Range("H112").Value = 79976.336
  $ws.Range("J112").Value = 79976.336
  $ws.Range("L112").Value = 79976.336
  $ws.Range("N112").Value = -82930.336
  # Row 113
  $ws.Range("H113").Value = 4883.25
  $ws.Range("J113").Value = 5673.625
  $ws.Range("L113").Value = 5673.625
  $ws.Range("N113").Value = -10013.625
  # Row 132
  $ws.Range("H132").Value = 3575.0417
  $ws.Range("I132").Value = 3800.15
  $ws.Range("K132").Value = 11400.45
  $ws.Range("M132").Value = -8870.450000000001
  # Row 134
  $ws.Range("H134").Value = 5442.8125
  $ws.Range("I134").Value = 4668.115
  $ws.Range("K134").Value = 14004.345
  $ws.Range("M134").Value = -11469.345
  # Row 136
  $ws.Range("H136").Value = 5273.7144
  $ws.Range("I136").Value = 4072.5
  $ws.Range("K136").Value = 12217.5
  $ws.Range("M136").Value = -9667.5
  # Row 141
  $ws.Range("H141").Value = 38800
  $ws.Range("J141").Value = 38800
  $ws.Range("L141").Value = 38800
  $ws.Range("N141").Value = -49160

$ws = $wb.Worksheets.Item("CUL")
  # Row 12
  $ws.Range("H12").Value = 338.8889
  $ws.Range("I12").Value = 3.5714285
  $ws.Range("J12").Value = 1512.5
  $ws.Range("K12").Value = 10.7142855
  $ws.Range("L12").Value = 4537.5
  $ws.Range("M12").Value = 162.2857145
  $ws.Range("N12").Value = -4883.5
  # Row 136
  $ws.Range("H136").Value = 0
  $ws.Range("I136").Value = 0
  $ws.Range("K136").Value = 0
  $ws.Range("M136").Value = $null
  # Row 138
  $ws.Range("H138").Value = 2047.7142
  $ws.Range("I138").Value = 2047.7142
  $ws.Range("K138").Value = 6143.142599999999
  $ws.Range("M138").Value = -1003.142599999999
  # Row 141
  $ws.Range("H141").Value = 1762.8889
  $ws.Range("I141").Value = 1762.8889
  $ws.Range("K141").Value = 5288.6667
  $ws.Range("M141").Value = -108.6666999999998

$ws = $wb.Worksheets.Item("GSM")
  # Row 63
  $ws.Range("H63").Value = 40800
  $ws.Range("J63").Value = 40800
  $ws.Range("L63").Value = 40800
  $ws.Range("N63").Value = -42172
  # Row 66
  $ws.Range("H66").Value = 40800
  $ws.Range("J66").Value = 40800
  $ws.Range("L66").Value = 122400
  $ws.Range("N66").Value = -129264
  # Row 70
  $ws.Range("H70").Value = 12979.667
  $ws.Range("I70").Value = 10216.75
  $ws.Range("J70").Value = 15190
  $ws.Range("K70").Value = 10216.75
  $ws.Range("L70").Value = 15190
  $ws.Range("M70").Value = -9946.75
  $ws.Range("N70").Value = -15730
  # Row 73
  $ws.Range("H73").Value = 12979.667
  $ws.Range("I73").Value = 10216.75
  $ws.Range("J73").Value = 15190
  $ws.Range("K73").Value = 10216.75
  $ws.Range("L73").Value = 15190
  $ws.Range("M73").Value = -9280.75
  $ws.Range("N73").Value = -17062
  # Row 102
  $ws.Range("H102").Value = 4289.1577
  $ws.Range("I102").Value = 2170.5833
  $ws.Range("K102").Value = 2170.5833
  $ws.Range("M102").Value = -548.5832999999998
  # Row 126
  $ws.Range("H126").Value = 6374.875
  $ws.Range("I126").Value = 5166.6665
  $ws.Range("K126").Value = 15499.9995
  $ws.Range("M126").Value = -13029.9995
  # Row 132
  $ws.Range("H132").Value = 1583.3448
  $ws.Range("I132").Value = 1457.8077
  $ws.Range("K132").Value = 4373.4231
  $ws.Range("M132").Value = -1843.4231

$ws = $wb.Worksheets.Item("LTW")
  # Row 20
  $ws.Range("H20").Value = 98333
  $ws.Range("I20").Value = 15000
  $ws.Range("K20").Value = 15000
  $ws.Range("M20").Value = -14774
  # Row 61
  $ws.Range("H61").Value = 84963.5
  $ws.Range("I61").Value = 101656.4
  $ws.Range("K61").Value = 101656.4
  $ws.Range("M61").Value = -101454.4
  # Row 113
  $ws.Range("H113").Value = 84963.5
  $ws.Range("I113").Value = 101656.4
  $ws.Range("K113").Value = 101656.4
  $ws.Range("M113").Value = -99486.39999999999
  # Row 132
  $ws.Range("H132").Value = 10836.702
  $ws.Range("I132").Value = 11679.771
  $ws.Range("K132").Value = 35039.313
  $ws.Range("M132").Value = -32509.313
  # Row 136
  $ws.Range("H136").Value = 4150.086
  $ws.Range("I136").Value = 3795.879
  $ws.Range("K136").Value = 11387.637
  $ws.Range("M136").Value = -8837.636999999999

$ws = $wb.Worksheets.Item("WVR")
  # Row 100
  $ws.Range("H100").Value = 450
  $ws.Range("I100").Value = 450
  $ws.Range("J100").Value = 0
  $ws.Range("K100").Value = 900
  $ws.Range("L100").Value = 0
  $ws.Range("M100").Value = -359
  $ws.Range("N100").Value = $null
  # Row 132
  $ws.Range("H132").Value = 5170.1816
  $ws.Range("I132").Value = 5189.8
  $ws.Range("K132").Value = 15569.4
  $ws.Range("M132").Value = -13039.4
